$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.906.13"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "2.221.19"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.27"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.49"
$ws.Range("E6").Value = "  +7.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.517"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.27"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.50"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.38"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "2.560.41"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.07"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "2.224.69"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.730"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").Value = "39.839.31"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.52"
$ws.Range("E20").Value = "  +11.85%  "
$ws.Range("D21").Value = "0.0₃0883"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.83"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.77"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.80"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.83"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.26"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.81"
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.06"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0720"
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.02"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0989"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").Value = "2.095.19"
$ws.Range("E42").Value = "  +9.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.79"
$ws.Range("E43").Value = "  +4.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  +5.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0270"
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.02"
$ws.Range("E46").Value = "  +9.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.75"
$ws.Range("E47").Value = "  +8.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "2.430.96"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.92"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.28"
$ws.Range("E51").Value = "  +1.28%  "
